$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.191.13"
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").Value = "3.379.45"
$ws.Range("E3").Value = "  +1.31%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.66"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.94"
$ws.Range("E6").Value = "  +1.35%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("E8").Value = "  +1.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("E9").Value = "  +5.64%  "

$ws.Range("E10").Value = "  +1.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.36"
$ws.Range("E11").Value = "  +2.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("E12").Value = "  +2.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "683.58"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("E14").Value = "  +1.89%  "

$ws.Range("D15").Value = "3.921.79"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").Value = "69.226.99"
$ws.Range("E16").Value = "  +2.06%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +1.55%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.385.02"
$ws.Range("E18").Value = "  +1.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.63"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.28"
$ws.Range("E20").Value = "  +1.85%  "

$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.41"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.15"
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("E24").Value = "  +3.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.71"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").Value = "  +0.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.98"
$ws.Range("E28").Value = "  +2.88%  "

$ws.Range("E29").Value = "  +1.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "559.09"
$ws.Range("E31").Value = "  -1.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.11"
$ws.Range("E32").Value = "  +0.86%  "

$ws.Range("E33").Value = "  +0.31%  "

$ws.Range("E34").Value = "  +5.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.76"
$ws.Range("E35").Value = "  +2.41%  "

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "3.677.57"
$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.64"
$ws.Range("E38").Value = "  +2.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("E39").Value = "  +3.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.27"
$ws.Range("E40").Value = "  +3.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("D42").Value = "0.0₃0696"
$ws.Range("E42").Value = "  +2.89%  "

$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0422"
$ws.Range("E44").Value = "  +3.43%  "

$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("E46").Value = "  -0.29%  "

$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.41"
$ws.Range("E48").Value = "  +4.70%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.21"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  +4.13%  "
